$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new date header "06-sep" in BK1
$ws.Range("BK1").Value = "06-sep"

# New values for the "06-sep" column (BK), one per data row 2-11
$values = @(16, 14, 11, 15, 8, 19, 22, 16, 17, 15)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 63).Value = $values[$i]
}

# Match the number format/alignment already used by the adjacent BJ column
$ws.Range("BJ2:BJ11").Copy()
$ws.Range("BK2:BK11").PasteSpecial(-4122)

# Update the active selection to match the post-edit state
$ws.Range("BK12").Select()
